$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 95634
$ws.Range("B2").Value = "Dr. Joaquim Sousa"
$ws.Range("C2").Value = "Juridico"
$ws.Range("F2").Value = 45090
$ws.Range("G2").Value = 2569.77

# Row 3
$ws.Range("A3").Value = 72897
$ws.Range("B3").Value = "Sra. Isabelly Pinto"
$ws.Range("C3").Value = "Juridico"
$ws.Range("D3").Value = "Outros"
$ws.Range("E3").Value = 5
$ws.Range("F3").Value = 45089
$ws.Range("G3").Value = 5923.98

# Row 4
$ws.Range("A4").Value = 83000
$ws.Range("B4").Value = "Ísis Oliveira"
$ws.Range("C4").Value = "TI"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 45083
$ws.Range("G4").Value = 6730.48

# Row 5
$ws.Range("A5").Value = 96119
$ws.Range("B5").Value = "Dom da Paz"
$ws.Range("C5").Value = "Engenharia"
$ws.Range("D5").Value = "Viagem de negocios"
$ws.Range("E5").Value = 6
$ws.Range("F5").Value = 45098
$ws.Range("G5").Value = 3345.14

# Row 6
$ws.Range("A6").Value = 6173
$ws.Range("B6").Value = "Sr. Felipe Cunha"
$ws.Range("C6").Value = "Operacoes"
$ws.Range("F6").Value = 45091
$ws.Range("G6").Value = 3863.1

# Row 7
$ws.Range("A7").Value = 18305
$ws.Range("B7").Value = "Thomas Mendes"
$ws.Range("C7").Value = "Marketing"
$ws.Range("E7").Value = 5
$ws.Range("F7").Value = 45100
$ws.Range("G7").Value = 5839.93

# Row 8
$ws.Range("A8").Value = 62093
$ws.Range("B8").Value = "Luiz Fernando Rezende"
$ws.Range("C8").Value = "Marketing"
$ws.Range("D8").Value = "Doenca"
$ws.Range("E8").Value = 5
$ws.Range("G8").Value = 6200.79

# Row 9
$ws.Range("A9").Value = 15442
$ws.Range("B9").Value = "Kamilly Gonçalves"
$ws.Range("D9").Value = "Viagem de negocios"
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 45103
$ws.Range("G9").Value = 7661.2

# Row 10
$ws.Range("A10").Value = 94542
$ws.Range("B10").Value = "José Miguel Alves"
$ws.Range("C10").Value = "Operacoes"
$ws.Range("D10").Value = "Problemas pessoais"
$ws.Range("E10").Value = 7
$ws.Range("F10").Value = 45086
$ws.Range("G10").Value = 2041.94

# Row 11
$ws.Range("A11").Value = 34619
$ws.Range("B11").Value = "Ana Liz Rezende"
$ws.Range("C11").Value = "Atendimento ao Cliente"
$ws.Range("D11").Value = "Outros"
$ws.Range("F11").Value = 45104
$ws.Range("G11").Value = 8266.67
